# Scheduled data refresh: updates currentAveragePrice* / Leve profit columns (H:N)
# in each job sheet with newly retrieved market-board figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41: The Write Stuff (Enchanted Mythril Ink)
$ws.Range("H41").Value = 372.76923
$ws.Range("I41").Value = 198
$ws.Range("K41").Value = 198
$ws.Range("M41").Value = 242

# Row 62: The Mustache Suits Him (Enchanted Mythrite Ink)
$ws.Range("H62").Value = 134626130
$ws.Range("I62").Value = 50013970
$ws.Range("J62").Value = 416666660
$ws.Range("K62").Value = 50013970
$ws.Range("L62").Value = 416666660
$ws.Range("M62").Value = -50013346
$ws.Range("N62").Value = -416667908

# Row 65: Forgery of Convenience (L) (Enchanted Mythrite Ink)
$ws.Range("H65").Value = 134626130
$ws.Range("I65").Value = 50013970
$ws.Range("J65").Value = 416666660
$ws.Range("K65").Value = 250069850
$ws.Range("L65").Value = 2083333300
$ws.Range("M65").Value = -250066730
$ws.Range("N65").Value = -2083339540

# Row 125: Body over Mind (Grade 5 Dexterity Alkahest)
$ws.Range("H125").Value = 1741596.1
$ws.Range("I125").Value = 2088
$ws.Range("J125").Value = 2437399.5
$ws.Range("K125").Value = 18792
$ws.Range("L125").Value = 21936595.5
$ws.Range("M125").Value = -16332
$ws.Range("N125").Value = -21941515.5

# Row 137: Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws.Range("H137").Value = 21370388
$ws.Range("I137").Value = 5209296.5
$ws.Range("J137").Value = 69853656
$ws.Range("K137").Value = 15627889.5
$ws.Range("L137").Value = 209560968
$ws.Range("M137").Value = -15625339.5
$ws.Range("N137").Value = -209566068

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff (Cobalt Ingot)
$ws.Range("H61").Value = 5912162.5
$ws.Range("I61").Value = 2778782.5
$ws.Range("J61").Value = 29412514
$ws.Range("K61").Value = 2778782.5
$ws.Range("L61").Value = 29412514
$ws.Range("M61").Value = -2778570.5
$ws.Range("N61").Value = -29412938

# Row 74: As the Bolt Flies (Titanium Nugget)
$ws.Range("H74").Value = 44976080
$ws.Range("I74").Value = 56350400
$ws.Range("J74").Value = 22227436
$ws.Range("K74").Value = 56350400
$ws.Range("L74").Value = 22227436
$ws.Range("M74").Value = -56349526
$ws.Range("N74").Value = -22229184

# Row 77: Heavy Metal Banned (L) (Titanium Nugget)
$ws.Range("H77").Value = 44976080
$ws.Range("I77").Value = 56350400
$ws.Range("J77").Value = 22227436
$ws.Range("K77").Value = 281752000
$ws.Range("L77").Value = 111137180
$ws.Range("M77").Value = -281747632
$ws.Range("N77").Value = -111145916

# Row 102: Smells of Rich Tama-hagane (Tama-hagane Ingot)
$ws.Range("H102").Value = 1945.5625
$ws.Range("I102").Value = 1966.3572
$ws.Range("J102").Value = 1800
$ws.Range("K102").Value = 1966.3572
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = -344.3571999999999
$ws.Range("N102").Value = -5044

# Row 122: Haste for High Durium (High Durium Nugget)
$ws.Range("H122").Value = 3552.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3552.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 10657.5
$ws.Range("N122").Value = -15557.5
$ws.Range("M122").ClearContents()

# Row 136: Metal with Mettle (Cobalt Tungsten Ingot)
$ws.Range("H136").Value = 5912162.5
$ws.Range("I136").Value = 2778782.5
$ws.Range("J136").Value = 29412514
$ws.Range("K136").Value = 8336347.5
$ws.Range("L136").Value = 88237542
$ws.Range("M136").Value = -8333797.5
$ws.Range("N136").Value = -88242642

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin (Adamantite Nugget)
$ws.Range("H86").Value = 1899.19
$ws.Range("I86").Value = 1899.19
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1899.19
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -776.1900000000001
$ws.Range("N86").ClearContents()

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) (Adamantite Nugget)
$ws.Range("H89").Value = 1899.19
$ws.Range("I89").Value = 1899.19
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9495.950000000001
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -3879.950000000001
$ws.Range("N89").ClearContents()

# Row 114: Halfhearted Effort (Bluespirit Halfheart Saw)
$ws.Range("H114").Value = 39750
$ws.Range("J114").Value = 39750
$ws.Range("L114").Value = 39750
$ws.Range("N114").Value = -48428

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found (Walnut Lumber)
$ws.Range("H31").Value = 2128602
$ws.Range("I31").Value = 1097863.8
$ws.Range("J31").Value = 5689334.5
$ws.Range("K31").Value = 1097863.8
$ws.Range("L31").Value = 5689334.5
$ws.Range("M31").Value = -1097568.8
$ws.Range("N31").Value = -5689924.5

# Row 34: Armoires of the Rich and Famous (Walnut Lumber)
$ws.Range("H34").Value = 2128602
$ws.Range("I34").Value = 1097863.8
$ws.Range("J34").Value = 5689334.5
$ws.Range("K34").Value = 1097863.8
$ws.Range("L34").Value = 5689334.5
$ws.Range("M34").Value = -1097661.8
$ws.Range("N34").Value = -5689738.5

# Row 87: Anatomy of a Drill Bit (Dragonscale Grinding Wheel)
$ws.Range("H87").Value = 31500
$ws.Range("J87").Value = 31500
$ws.Range("L87").Value = 31500
$ws.Range("N87").Value = -33872

# Row 90: Pulling Them to the Grind (L) (Dragonscale Grinding Wheel)
$ws.Range("H90").Value = 31500
$ws.Range("J90").Value = 31500
$ws.Range("L90").Value = 94500
$ws.Range("N90").Value = -106356

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap (Maple Syrup)
$ws.Range("H5").Value = 2992266
$ws.Range("I5").Value = 2137458.5
$ws.Range("K5").Value = 6412375.5
$ws.Range("M5").Value = -6412263.5

# Row 54: Good Eats in Ishgard (Salt Cod Puffs)
$ws.Range("H54").Value = 1840
$ws.Range("I54").Value = 1200
$ws.Range("K54").Value = 3600
$ws.Range("M54").Value = -3041

# Row 122: Salt of the North (Northern Sea Salt)
$ws.Range("H122").Value = 800.6
$ws.Range("I122").Value = 377.7
$ws.Range("K122").Value = 3399.3
$ws.Range("M122").Value = -949.2999999999997

# Row 132: More Mezcal (Cooking Mezcal)
$ws.Range("H132").Value = 1843.0625
$ws.Range("I132").Value = 1838
$ws.Range("J132").Value = 1848.125
$ws.Range("K132").Value = 16542
$ws.Range("L132").Value = 16633.125
$ws.Range("M132").Value = -14012
$ws.Range("N132").Value = -21693.125

# Row 135: Not-so-secret Ingredient (Royal Maple Syrup)
$ws.Range("H135").Value = 2992266
$ws.Range("I135").Value = 2137458.5
$ws.Range("K135").Value = 19237126.5
$ws.Range("M135").Value = -19234591.5

$ws = $wb.Worksheets.Item("GSM")
# Row 11: A Ringing Success (Copper Ring)
$ws.Range("H11").Value = 84002000
$ws.Range("I11").Value = 105001250
$ws.Range("K11").Value = 105001250
$ws.Range("M11").Value = -105001111

# Row 69: High Above Me, She Sews Lovely (Mythrite Needle)
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 72: Old-school Spooling (L) (Mythrite Needle)
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Row 122: Awarding Academic Excellence (Ametrine)
$ws.Range("H122").Value = 20836248
$ws.Range("I122").Value = 3521.4
$ws.Range("J122").Value = 55557456
$ws.Range("K122").Value = 10564.2
$ws.Range("L122").Value = 166672368
$ws.Range("M122").Value = -8114.200000000001
$ws.Range("N122").Value = -166677268

$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack (Tiger Leather)
$ws.Range("H100").Value = 1519.6471
$ws.Range("I100").Value = 1376.1538
$ws.Range("J100").Value = 1986
$ws.Range("K100").Value = 1376.1538
$ws.Range("L100").Value = 1986
$ws.Range("M100").Value = -835.1538
$ws.Range("N100").Value = -3068

# Row 122: Hell on Leather (Gaja Leather)
$ws.Range("H122").Value = 18442556
$ws.Range("I122").Value = 2366408.8
$ws.Range("J122").Value = 66671000
$ws.Range("K122").Value = 7099226.399999999
$ws.Range("L122").Value = 200013000
$ws.Range("M122").Value = -7096776.399999999
$ws.Range("N122").Value = -200017900

$ws = $wb.Worksheets.Item("WVR")
# Row 100: Of Great Import (Kudzu Thread)
$ws.Range("H100").Value = 17200
$ws.Range("I100").Value = 20560
$ws.Range("K100").Value = 41120
$ws.Range("M100").Value = -40579
